$d = $word.ActiveDocument

# The "Uddannelsesmæssig Baggrund" (education history) table is the only
# table in the document. Row 1 is the header (Skoleforløb / Fra / Til) and
# row 2 is "Hovedforløb 3". A new row for "Skolepraktik IV" needs to be
# inserted right after the header row, i.e. immediately before "Hovedforløb 3".

$t = $d.Tables.Item(1)
$followingRow = $t.Rows.Item(2)

$newRow = $t.Rows.Add($followingRow)

$newRow.Cells.Item(1).Range.Text = "Skolepraktik IV"
$newRow.Cells.Item(2).Range.Text = "21. jun. 2021"
$newRow.Cells.Item(3).Range.Text = "Nu"
